$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new task row (row 16): "15. Add button check all checkboxes"
$ws.Range("A16").Value = "15. Add button check all checkboxes"
$ws.Range("B16").Value = "Low"
$ws.Range("C16").Value = "Open"

# Highlight the "In progress" status cell (C5) in red font
$ws.Range("C5").Font.Color = 255

# Move the active selection to C16, matching the saved view state
[void]$ws.Range("C16").Select()
